$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: updated TPM-derived statistics (columns A-H and the
#     Sending/Ligand/Receptor/Target cluster labels are unchanged) ---
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04936366666666667
$ws.Range("N2").Value = 0.148091
$ws.Range("O2").Value = 0.04616170608573571
$ws.Range("P2").Value = 0.0461617060857357
$ws.Range("Q2").Value = 0.002366773907444444
$ws.Range("R2").Value = 0.021300965167
$ws.Range("S2").Value = 0.04616170608573571
$ws.Range("T2").Value = 0.0461617060857357

# --- Row 3 ---
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.008057748967298944
$ws.Range("P3").Value = 0.008057748967298944
$ws.Range("S3").Value = 0.008057748967298944
$ws.Range("T3").Value = 0.008057748967298944

# --- Row 4 ---
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 1.011383666666666
$ws.Range("N4").Value = 3.034151
$ws.Range("O4").Value = 0.9457805449469654
$ws.Range("P4").Value = 0.9457805449469653
$ws.Range("Q4").Value = 0.0484914641541111
$ws.Range("R4").Value = 0.4364231773869999
$ws.Range("S4").Value = 0.9457805449469654
$ws.Range("T4").Value = 0.9457805449469653

# --- Rows 5-7 no longer exist in the refreshed (TPM) output; the old
#     MuSCs-target rows were folded away, leaving just the 3 data rows
#     above for the FAPs sending cluster. ---
$ws.Rows("5:7").Delete()
